# Update prediction/stats data for Bayer 04 Leverkusen player stats sheet
# (sofaplayer/Bundesliga/Bayer 04 Leverkusen_stats.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Christian Kofane
$ws.Range("E2").Value = 6.7222222222222
$ws.Range("K2").Value = 3.7862
$ws.Range("AA2").Value = 0.42111274
$ws.Range("AB2").Value = 231
$ws.Range("AE2").Value = 95
$ws.Range("AF2").Value = 84.07079646017699
$ws.Range("AI2").Value = 65
$ws.Range("AJ2").Value = 44
$ws.Range("BG2").Value = 51
$ws.Range("BQ2").Value = 121
$ws.Range("BT2").Value = 18

# Row 3: Patrik Schick
$ws.Range("E3").Value = 6.8705882352941
$ws.Range("AF3").Value = 79.766536964981
$ws.Range("AG3").Value = 257
$ws.Range("AP3").Value = 4
$ws.Range("BA3").Value = 57
$ws.Range("BB3").Value = 53.271028037383
$ws.Range("BC3").Value = 18
$ws.Range("BD3").Value = 39.130434782609
$ws.Range("BG3").Value = 113
$ws.Range("BQ3").Value = 116.8
$ws.Range("BT3").Value = 52
$ws.Range("DB3").Value = 50
$ws.Range("DF3").Value = 99
$ws.Range("DG3").Value = 158

# Row 4: Jonas Hofmann
$ws.Range("AA4").Value = 1.61113524

# Row 5: Martin Terrier
$ws.Range("AB5").Value = 236

# Row 7: Ibrahim Maza
$ws.Range("AA7").Value = 1.36198075
$ws.Range("AB7").Value = 712

# Row 9: Alejandro Grimaldo
$ws.Range("K9").Value = 4.0733
$ws.Range("AA9").Value = 3.744844
$ws.Range("AE9").Value = 963
$ws.Range("AF9").Value = 87.78486782133101
$ws.Range("AG9").Value = 1097
$ws.Range("AI9").Value = 542
$ws.Range("AR9").Value = 73
$ws.Range("AZ9").Value = 35.555555555556
$ws.Range("CB9").Value = 11
$ws.Range("CS9").Value = 45
$ws.Range("DG9").Value = 628

# Row 10: Exequiel Palacios
$ws.Range("AA10").Value = 0.16340783
$ws.Range("AE10").Value = 208
$ws.Range("AF10").Value = 90.434782608696
$ws.Range("AI10").Value = 127
$ws.Range("AJ10").Value = 40
$ws.Range("BG10").Value = 29
$ws.Range("BT10").Value = 22
$ws.Range("CD10").Value = 3

# Row 12: Malik Tillman
$ws.Range("E12").Value = 7.0411764705882
$ws.Range("K12").Value = 2.6507
$ws.Range("BQ12").Value = 119.7

# Row 13: Ernest Poku
$ws.Range("I13").Value = 2
$ws.Range("K13").Value = 2.7357
$ws.Range("AF13").Value = 84.375
$ws.Range("AG13").Value = 416
$ws.Range("BT13").Value = 65
$ws.Range("DG13").Value = 312

# Row 14: Aleix Garcia
$ws.Range("AA14").Value = 4.00135481
$ws.Range("AB14").Value = 2357
$ws.Range("AE14").Value = 1969
$ws.Range("AF14").Value = 93.406072106262
$ws.Range("AG14").Value = 2108
$ws.Range("AI14").Value = 1050
$ws.Range("AJ14").Value = 438
$ws.Range("DG14").Value = 1143

# Row 15: Robert Andrich
$ws.Range("AA15").Value = 0.98787781
$ws.Range("AB15").Value = 1316
$ws.Range("AE15").Value = 1018
$ws.Range("AF15").Value = 89.84995586937301
$ws.Range("AG15").Value = 1133
$ws.Range("AH15").Value = 564
$ws.Range("CC15").Value = 64
$ws.Range("DF15").Value = 599

# Row 17: Lucas Vazquez
$ws.Range("E17").Value = 7.0125
$ws.Range("AA17").Value = 1.05735107
$ws.Range("AE17").Value = 130
$ws.Range("AF17").Value = 86.09271523178801
$ws.Range("AH17").Value = 44
$ws.Range("AL17").Value = 33.333333333333
$ws.Range("BB17").Value = 56.25
$ws.Range("BD17").Value = 53.846153846154
$ws.Range("BQ17").Value = 56.1
$ws.Range("BT17").Value = 21
$ws.Range("CD17").Value = 3
$ws.Range("CU17").Value = 14
$ws.Range("CZ17").Value = 6
$ws.Range("DF17").Value = 51
$ws.Range("DG17").Value = 100

# Row 18: Edmond Tapsoba
$ws.Range("I18").Value = 3
$ws.Range("AA18").Value = 0.67841962
$ws.Range("BB18").Value = 61.702127659574
$ws.Range("BF18").Value = 63.636363636364
$ws.Range("CC18").Value = 66
$ws.Range("CD18").Value = 37
$ws.Range("CU18").Value = 36
$ws.Range("CV18").Value = 20

# Row 19: Jarell Quansah
$ws.Range("K19").Value = 0.865
$ws.Range("AA19").Value = 0.76336699
$ws.Range("AB19").Value = 1325
$ws.Range("AP19").Value = 17
$ws.Range("BA19").Value = 67
$ws.Range("BB19").Value = 60.909090909091
$ws.Range("BC19").Value = 44
$ws.Range("BD19").Value = 61.971830985915
$ws.Range("BG19").Value = 150
$ws.Range("DA19").Value = 14
$ws.Range("DB19").Value = 82.35294117647101

# Row 21: Arthur
$ws.Range("AA21").Value = 1.35296295

# Row 24: Janis Blaswich
$ws.Range("AE24").Value = 95
$ws.Range("AF24").Value = 81.19658119658099
$ws.Range("AG24").Value = 117
$ws.Range("AH24").Value = 88
$ws.Range("DF24").Value = 95
